# "every investment the capacity market price is recalculated" /
# "capacity market was overinvinvesting"
#
# Coupling Parameters sheet:
#  - End Year (B4): 2055 -> 2060
#  - start_dismantling_tick (B23): 5 -> 1
#  - scenarioWeatheryearsExcel (B32): switch weather-year file from the
#    "HalfFlexElectrolyzers" variant to the plain 40weatherYears2050TNO.xlsx
#    file, keeping the old file name around in the next column (D32) for
#    reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# End Year
$ws.Range("B4").Value2 = 2060

# start_dismantling_tick
$ws.Range("B23").Value2 = 1

# scenarioWeatheryearsExcel: keep a record of the previous file name in D32,
# then point B32 at the new file.
$ws.Range("D32").Value2 = $ws.Range("B32").Value2
$ws.Range("B32").Value2 = "40weatherYears2050TNO.xlsx"

# Restore the selection the author left the sheet on.
$ws.Range("F21").Select()
